# Auto-generated edit script: apply numeric corrections to the Leve profit
# calculation columns (H..N) across all eight job sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 195.3
$ws.Range("J53").Value = 180.875
$ws.Range("L53").Value = 180.875
$ws.Range("N53").Value = -1454.875
$ws.Range("H112").Value = 2829.261
$ws.Range("I112").Value = 1709.1428
$ws.Range("J112").Value = 3319.3125
$ws.Range("K112").Value = 5127.428400000001
$ws.Range("L112").Value = 9957.9375
$ws.Range("M112").Value = -4019.428400000001
$ws.Range("N112").Value = -12173.9375

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 808216.5600000001
$ws.Range("I32").Value = 869332.1
$ws.Range("K32").Value = 869332.1
$ws.Range("M32").Value = -869045.1
$ws.Range("H45").Value = 4377.846
$ws.Range("I45").Value = 3849.8572
$ws.Range("K45").Value = 3849.8572
$ws.Range("M45").Value = -3472.8572
$ws.Range("H61").Value = 2567805.2
$ws.Range("I61").Value = 3571.4375
$ws.Range("K61").Value = 3571.4375
$ws.Range("M61").Value = -3359.4375
$ws.Range("H132").Value = 7942.875
$ws.Range("I132").Value = 6044.9287
$ws.Range("K132").Value = 18134.7861
$ws.Range("M132").Value = -15604.7861
$ws.Range("H136").Value = 2567805.2
$ws.Range("I136").Value = 3571.4375
$ws.Range("K136").Value = 10714.3125
$ws.Range("M136").Value = -8164.3125

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 9999.3125
$ws.Range("J138").Value = 9999
$ws.Range("L138").Value = 9999
$ws.Range("N138").Value = -20279

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 5250
$ws.Range("J14").Value = 5500
$ws.Range("L14").Value = 5500
$ws.Range("N14").Value = -5840
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9826
$ws.Range("H31").Value = 1504960.4
$ws.Range("I31").Value = 2025348.6
$ws.Range("J31").Value = 1616.4445
$ws.Range("K31").Value = 2025348.6
$ws.Range("L31").Value = 1616.4445
$ws.Range("M31").Value = -2025053.6
$ws.Range("N31").Value = -2206.4445
$ws.Range("H34").Value = 1504960.4
$ws.Range("I34").Value = 2025348.6
$ws.Range("J34").Value = 1616.4445
$ws.Range("K34").Value = 2025348.6
$ws.Range("L34").Value = 1616.4445
$ws.Range("M34").Value = -2025146.6
$ws.Range("N34").Value = -2020.4445

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2732
$ws.Range("I17").Value = 2732
$ws.Range("K17").Value = 8196
$ws.Range("M17").Value = -8027
$ws.Range("H87").Value = 17887.889
$ws.Range("I87").Value = 10247.75
$ws.Range("K87").Value = 30743.25
$ws.Range("M87").Value = -29495.25
$ws.Range("H88").Value = 38599.4
$ws.Range("I88").Value = 44999.668
$ws.Range("J88").Value = 28999
$ws.Range("K88").Value = 134999.004
$ws.Range("L88").Value = 86997
$ws.Range("M88").Value = -134571.004
$ws.Range("N88").Value = -87853
$ws.Range("H90").Value = 17887.889
$ws.Range("I90").Value = 10247.75
$ws.Range("K90").Value = 92229.75
$ws.Range("M90").Value = -85989.75
$ws.Range("H91").Value = 38599.4
$ws.Range("I91").Value = 44999.668
$ws.Range("J91").Value = 28999
$ws.Range("K91").Value = 134999.004
$ws.Range("L91").Value = 86997
$ws.Range("M91").Value = -133517.004
$ws.Range("N91").Value = -89961
$ws.Range("H121").Value = 18833.572
$ws.Range("I121").Value = 370.75
$ws.Range("K121").Value = 1112.25
$ws.Range("M121").Value = 197.75
$ws.Range("H131").Value = 6965.25
$ws.Range("J131").Value = 6965.25
$ws.Range("L131").Value = 20895.75
$ws.Range("N131").Value = -30975.75
$ws.Range("H137").Value = 7376.467
$ws.Range("I137").Value = 2841.5
$ws.Range("J137").Value = 10399.777
$ws.Range("K137").Value = 8524.5
$ws.Range("L137").Value = 31199.331
$ws.Range("M137").Value = -3424.5
$ws.Range("N137").Value = -41399.331
$ws.Range("H140").Value = 2083
$ws.Range("I140").Value = 1624.5
$ws.Range("K140").Value = 4873.5
$ws.Range("M140").Value = 306.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5937.3335
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H92").Value = 7062.75
$ws.Range("J92").Value = 7062.75
$ws.Range("L92").Value = 7062.75
$ws.Range("N92").Value = -10806.75
$ws.Range("H132").Value = 28245
$ws.Range("I132").Value = 28245
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 84735
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -82205
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 44999.5
$ws.Range("J134").Value = 44999.5
$ws.Range("L134").Value = 134998.5
$ws.Range("N134").Value = -140068.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 4700.6924
$ws.Range("I46").Value = 1280.5
$ws.Range("K46").Value = 1280.5
$ws.Range("M46").Value = -1092.5
$ws.Range("H122").Value = 5098.6665
$ws.Range("I122").Value = 2650.5
$ws.Range("K122").Value = 7951.5
$ws.Range("M122").Value = -5501.5
$ws.Range("H132").Value = 5132273
$ws.Range("I132").Value = 6670154.5
$ws.Range("K132").Value = 20010463.5
$ws.Range("M132").Value = -20007933.5
$ws.Range("H136").Value = 6946828.5
$ws.Range("I136").Value = 4034494.8
$ws.Range("K136").Value = 12103484.4
$ws.Range("M136").Value = -12100934.4

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1696.25
$ws.Range("I23").Value = 265.83334
$ws.Range("J23").Value = 5987.5
$ws.Range("K23").Value = 265.83334
$ws.Range("L23").Value = 5987.5
$ws.Range("M23").Value = -36.83334000000002
$ws.Range("N23").Value = -6445.5
$ws.Range("H44").Value = 237
$ws.Range("I44").Value = 237
$ws.Range("K44").Value = 237
$ws.Range("M44").Value = 317
$ws.Range("H81").Value = 146785.86
$ws.Range("I81").Value = 4500
$ws.Range("J81").Value = 253500.25
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 507000.5
$ws.Range("M81").Value = -7939
$ws.Range("N81").Value = -509122.5
$ws.Range("H84").Value = 146785.86
$ws.Range("I84").Value = 4500
$ws.Range("J84").Value = 253500.25
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 2535002.5
$ws.Range("M84").Value = -39696
$ws.Range("N84").Value = -2545610.5
$ws.Range("H132").Value = 3877056
$ws.Range("I132").Value = 4630587.5
$ws.Range("J132").Value = 1750.5714
$ws.Range("K132").Value = 13891762.5
$ws.Range("L132").Value = 5251.7142
$ws.Range("M132").Value = -13889232.5
$ws.Range("N132").Value = -10311.7142
